# Adds the HUDM-006 "Historia de Usuario" block right after the
# HUDM-005 block's closing paragraph ("Para: Acceder rápidamente...").
#
# New paragraphs inserted (in order):
#   1. "HUDM-006 Editar un Proveedor para contacto"  (3 runs: "HUDM-00", "6",
#      " Editar un Proveedor para contacto" -- matching the HUDM-004/HUDM-005
#      pattern already used elsewhere in the doc)
#   2. "Como: Supervisor"
#   3. "Quiero: Editar el correo o teléfono de un proveedor existente."
#   4. "Para: Asegurar que los datos estén siempre actualizados y no se
#      pierda la comunicación de contacto."

$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# Locate the paragraph that ends the HUDM-005 story ("Para: Acceder
# rápidamente..."). Keep using the SAME range object for Find + the
# position read straight after, since re-deriving $d.Content later
# returns a brand-new whole-document range, not the found location.
$findRng = $d.Content
$anchorText = "Para: Acceder rápidamente a su información de contacto sin tener que revisar toda la base de datos."
$found = $findRng.Find.Execute($anchorText, $true, $false, $false, $false, $false, `
                                $true, 1, $false, "", 0)
if (-not $found) {
    throw "Anchor paragraph not found"
}
$anchorEnd = $findRng.End

# Split right after the anchor paragraph's text (before its own paragraph
# mark), producing a brand-new empty paragraph right after it.
$d.Range($anchorEnd, $anchorEnd).InsertParagraphAfter()
$curPara = $d.Range($anchorEnd + 1, $anchorEnd + 1).Paragraphs(1)

# Ordered (run-texts-per-paragraph) content to insert. The first new
# paragraph is split into three runs; the rest are single plain runs.
$paragraphsToAdd = @(
    , @('HUDM-00', '6', ' Editar un Proveedor para contacto')
    , @('Como: Supervisor')
    , @('Quiero: Editar el correo o teléfono de un proveedor existente.')
    , @('Para: Asegurar que los datos estén siempre actualizados y no se pierda la comunicación de contacto.')
)

for ($i = 0; $i -lt $paragraphsToAdd.Count; $i++) {
    $runTexts = $paragraphsToAdd[$i]

    # Replace the (currently empty) paragraph's full range -- start through
    # its own paragraph mark -- with our run markup, so no stray empty run
    # or borrowed pPr/rPr formatting from a neighbour paragraph survives.
    $runsXml = ($runTexts | ForEach-Object {
        $t = $_
        if ($t -match '^\s' -or $t -match '\s$') {
            '<w:r><w:t xml:space="preserve">' + $t + '</w:t></w:r>'
        } else {
            '<w:r><w:t>' + $t + '</w:t></w:r>'
        }
    }) -join ''
    $fragment = '<w:p ' + $wNs + '>' + $runsXml + '</w:p>'

    $target = $d.Range($curPara.Range.Start, $curPara.Range.End)
    $target.InsertXML($fragment)

    # Re-resolve the paragraph we just filled (Start is stable across the
    # InsertXML call) so we can read its real End (start of the following
    # paragraph) and chain the next insertion after it.
    $curPara = $d.Range($curPara.Range.Start, $curPara.Range.Start).Paragraphs(1)

    if ($i -lt ($paragraphsToAdd.Count - 1)) {
        $curPara.Range.InsertParagraphAfter()
        $curPara = $d.Range($curPara.Range.End, $curPara.Range.End).Paragraphs(1)
    }
}

Write-Output "Done. Paragraph count now: $($d.Paragraphs.Count)"
